$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("value") now stores the solver numeric result as TEXT
# (matching the source data), not as a native number, so force the
# number format to Text before writing so Excel keeps it as a string.
$ws.Range("C2:C12").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = '                              '
$ws.Range("B2").Value = [double]'0.0'
$ws.Range("C2").Value = '                              '
$ws.Range("D2").Value = '                                                                                                                                                                                                        '
$ws.Range("E2").Value = [double]'0.0'
$ws.Range("F2").Value = [double]'0.0'

# Row 3
$ws.Range("A3").Value = '10_ulysses_3.tsp'
$ws.Range("B3").Value = [double]'3.868'
$ws.Range("C3").Value = '136.99527629589417'
$ws.Range("D3").Value = '[5, 9][1, 2, 3, 10][4, 6, 7, 8]'
$ws.Range("E3").Value = [double]'136.99420807981383'
$ws.Range("F3").Value = [double]'7.797466520183073e-06'

# Row 4
$ws.Range("A4").Value = '10_ulysses_6.tsp'
$ws.Range("B4").Value = [double]'15.50600004196167'
$ws.Range("C4").Value = '55.11939124322688'
$ws.Range("D4").Value = '[10][6, 9][5][2, 3][1, 4, 8][7]'
$ws.Range("E4").Value = [double]'55.11442110891486'
$ws.Range("F4").Value = [double]'9.017034114322455e-05'

# Row 5
$ws.Range("A5").Value = '10_ulysses_9.tsp'
$ws.Range("B5").Value = [double]'7.161999940872192'
$ws.Range("C5").Value = '33.29189782877749'
$ws.Range("D5").Value = '[6][5][9][7][1][2, 3][8][4][10]'
$ws.Range("E5").Value = [double]'33.29189782877749'
$ws.Range("F5").Value = [double]'0.0'

# Row 6
$ws.Range("A6").Value = '14_burma_3.tsp'
$ws.Range("B6").Value = [double]'4.638000011444092'
$ws.Range("C6").Value = '93.38998725996821'
$ws.Range("D6").Value = '[3, 4, 5, 6, 7][1, 8, 9, 10, 11][2, 12, 13, 14]'
$ws.Range("E6").Value = [double]'93.38998725996821'
$ws.Range("F6").Value = [double]'0.0'

# Row 7
$ws.Range("A7").Value = '14_burma_6.tsp'
$ws.Range("B7").Value = [double]'30.003999948501587'
$ws.Range("C7").Value = '43.129462546693226'
$ws.Range("D7").Value = '[5, 10][7][2, 8][6, 12, 13][3, 4, 14][1, 9, 11]'
$ws.Range("E7").Value = [double]'0.0'
$ws.Range("F7").Value = [double]'0.9999999999976814'

# Row 8
$ws.Range("A8").Value = '14_burma_9.tsp'
$ws.Range("B8").Value = [double]'30.026999950408936'
$ws.Range("C8").Value = '20.762438566071065'
$ws.Range("D8").Value = '[3, 4][9, 11][10][1, 8][7][5][13, 14][6, 12][2]'
$ws.Range("E8").Value = [double]'0.0'
$ws.Range("F8").Value = [double]'0.9999999999951837'

# Row 9
$ws.Range("A9").Value = '22_ulysses_3.tsp'
$ws.Range("B9").Value = [double]'30.028000116348267'
$ws.Range("C9").Value = '515.5925785532522'
$ws.Range("D9").Value = '[1, 2, 3, 7, 13, 16, 18, 22][5, 6, 9, 10, 12, 14, 21][4, 8, 11, 15, 17, 19, 20]'
$ws.Range("E9").Value = [double]'149.9123235516313'
$ws.Range("F9").Value = [double]'0.7092426660361273'

# Row 10
$ws.Range("A10").Value = '26_eil_3.tsp'
$ws.Range("B10").Value = [double]'30.0239999294281'
$ws.Range("C10").Value = '3475.635504853873'
$ws.Range("D10").Value = '[6, 9, 13, 14, 16, 17, 18, 26][2, 4, 7, 8, 12, 20, 21, 22, 24][1, 3, 5, 10, 11, 15, 19, 23, 25]'
$ws.Range("E10").Value = [double]'191.3265443301569'
$ws.Range("F10").Value = [double]'0.9449520687474116'

# Row 11
$ws.Range("A11").Value = '26_eil_6.tsp'
$ws.Range("B11").Value = [double]'30.049000024795532'
$ws.Range("C11").Value = '1485.7858771055776'
$ws.Range("D11").Value = '[4, 13, 19][2, 5, 6, 10, 20][1, 14, 15, 16, 17, 25][3, 11, 12, 22, 23][7, 8, 24, 26][9, 18, 21]'
$ws.Range("E11").Value = [double]'0.0'
$ws.Range("F11").Value = [double]'0.9999999999999327'

# Row 12
$ws.Range("A12").Value = '26_eil_9.tsp'
$ws.Range("B12").Value = [double]'30.031999826431274'
$ws.Range("C12").Value = '1006.5109586516216'
$ws.Range("D12").Value = '[3, 8, 18][22, 23, 25][5, 14, 17][1, 4, 15][6, 12, 13][10, 11, 21][2, 7, 20][19, 24][9, 16, 26]'
$ws.Range("E12").Value = [double]'0.0'
$ws.Range("F12").Value = [double]'0.9999999999999006'
